# Refresh crypto price/volume snapshot data (GitHub Actions scheduled update).
# Column D mixes plain decimals ("7.41") with multi-dot thousands-grouped
# prices ("71.496.29"); Column E keeps literal leading/trailing padding
# spaces around the percentage text. Force the specific target cells to
# Text format first so Excel does not reinterpret/round them as numbers
# (e.g. "29.60" -> 29.6) while leaving every other cell style untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.496.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.805.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '701.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.19'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.802.37'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.45%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.41'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.46%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.75'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.451.75'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.807.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.495.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.73'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.57%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.01%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.09'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.06'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.716'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.59'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000143'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.33'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.55'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.55%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.958.13'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.13'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.28%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.58'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.61%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.60'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.32%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.37'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +12.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.44'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.04'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.994'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.88%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.20%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '164.54'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.53%  '

$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000308'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '44.89'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.59'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.303'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '418.07'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.67'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.90%  '

